# Update latest output (run 79)
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Schedule": shift the 3 schedule blocks up by one slot and
# drop the now-redundant 4th row (dimension shrinks from F4 to F3).
# ---------------------------------------------------------------
$schedule = $wb.Worksheets.Item("Schedule")

$schedule.Range("A2").Value2 = 46040.29166666666
$schedule.Range("B2").Value2 = 46040.79166666666
$schedule.Range("C2").Value2 = 12
$schedule.Range("D2").Value2 = 45.36
$schedule.Range("E2").Value2 = 213.75222375
$schedule.Range("F2").Value2 = 4.712350611772487

$schedule.Range("A3").Value2 = 46040.83333333334
$schedule.Range("B3").Value2 = 46041
$schedule.Range("C3").Value2 = 4
$schedule.Range("D3").Value2 = 15.12
$schedule.Range("E3").Value2 = 361.68226575
$schedule.Range("F3").Value2 = 23.9207847718254

# Row 4 no longer exists in the updated result - remove it entirely.
$schedule.Rows.Item(4).Delete()

# ---------------------------------------------------------------
# Sheet "Detailed": refresh the price/status columns with the
# latest run values.
# ---------------------------------------------------------------
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("E15").Value = "OFF"

$detailed.Range("B19").Value2 = 0.009549999999999999

$detailed.Range("B20").Value2 = 32.97732

$detailed.Range("B21").Value2 = 35.88
$detailed.Range("C21").Value = "historical"

$detailed.Range("B22").Value2 = 36.06
$detailed.Range("C22").Value = "historical"

$detailed.Range("B23").Value2 = 34.45798
$detailed.Range("E23").Value = "ON"

$detailed.Range("B25").Value2 = 0.7

$detailed.Range("B26").Value2 = 25.09786

$detailed.Range("B27").Value2 = 33.41567

$detailed.Range("B28").Value2 = 0.51

$detailed.Range("B29").Value2 = 0.01819

$detailed.Range("B30").Value2 = -0.32596

$detailed.Range("B31").Value2 = -4.22974

$detailed.Range("B32").Value2 = -4.23274

$detailed.Range("B33").Value2 = -7.48187

$detailed.Range("B34").Value2 = -2.81401

$detailed.Range("B35").Value2 = -5.8994

$detailed.Range("B36").Value2 = -5.79576

$detailed.Range("B37").Value2 = 5.18238

$detailed.Range("B38").Value2 = 6.98736

$detailed.Range("B39").Value2 = 18.62153

$detailed.Range("B41").Value2 = 55.05565

$detailed.Range("B42").Value2 = 56.72399

$detailed.Range("B43").Value2 = 46.34977

$detailed.Range("B44").Value2 = 56.98

$detailed.Range("B45").Value2 = 45.60231

$detailed.Range("B46").Value2 = 56.98

$detailed.Range("B47").Value2 = 36.2
